# wrapped up the BFS problems.
#
# Row 24 (2025-09-21, "Binary Tree - BFS - LC75+") gets its Easy/Medium
# quantity corrected and the Leetcode Problem # / Confidence Level filled
# in now that the BFS problems are wrapped up.
#
# Row 25 (2025-09-22) is renamed from "Binary Search Tree - LC75+" to
# "Binary Search Tree - LC75" (the placeholder "+" tier is dropped now
# that the topic has actually been started) and its Easy Quantity is
# corrected down from 3 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2025")

# Row 24 - finish filling in the BFS wrap-up row
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = "1161, 102"
$ws.Range("F24").Value = "High"

# Row 25 - rename the algorithm and correct the easy-quantity count
$ws.Range("B25").Value = "Binary Search Tree - LC75"
$ws.Range("C25").Value = 1

# Move the active selection to reflect where editing left off
$ws.Activate()
$ws.Range("F24").Select()
